$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.070.64"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "2.577.91"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "3.037.25"
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "58.023.44"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "2.576.34"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.421"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "0.0₃0725"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "36.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.830"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.816"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "280.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.89%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.588"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0948"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0532"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0227"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").Value = "1.912.46"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.92%  "
